$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = [DateTime]::FromOADate(45180)
$lastRow = 321

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value()
    if ($val -eq $oldDate -or $val -eq 45180) {
        $cell.Value = 45181
    }
}
